$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (Ajo / Chino / Primera, Feria Lagunitas de Puerto
# Montt) needs to be inserted as row 155; every existing record from the old
# row 155 onward shifts down by one row (old row 189 becomes row 190).
$ws.Rows.Item(155).Insert()

$ws.Range("A155").Value = 4
$ws.Range("B155").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C155").Value = "Los Lagos"
$ws.Range("D155").Value = 44543
$ws.Range("E155").Value = 10
$ws.Range("F155").Value = 100112003
$ws.Range("G155").Value = "Ajo"
$ws.Range("H155").Value = "Chino"
$ws.Range("I155").Value = "Primera"
$ws.Range("J155").Value = 40
$ws.Range("K155").Value = 21000
$ws.Range("L155").Value = 22000
$ws.Range("M155").Value = 21500
$ws.Range("N155").Value = "$/caja 10 kilos"
$ws.Range("O155").Value = "China"
$ws.Range("P155").Value = 2150
$ws.Range("Q155").Value = 10
$ws.Range("R155").Value = "Hortaliza"
